$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.356.21"
$ws.Range("E2").Value = "  +1.10%  "

$ws.Range("D3").Value = "3.357.12"
$ws.Range("E3").Value = "  +0.76%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.32"
$ws.Range("E5").Value = "  +0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.66"
$ws.Range("E6").Value = "  +0.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  +0.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.184"
$ws.Range("E9").Value = "  +3.78%  "

$ws.Range("E10").Value = "  +0.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.02"
$ws.Range("E11").Value = "  +5.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000273"
$ws.Range("E12").Value = "  +1.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "689.45"
$ws.Range("E13").Value = "  +2.68%  "

$ws.Range("D14").Value = "3.904.88"
$ws.Range("E14").Value = "  +0.63%  "

$ws.Range("E15").Value = "  +0.59%  "

$ws.Range("D16").Value = "68.371.58"
$ws.Range("E16").Value = "  +1.01%  "

$ws.Range("E17").Value = "  +1.32%  "

$ws.Range("D18").Value = "3.342.47"
$ws.Range("E18").Value = "  +0.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.49"
$ws.Range("E19").Value = "  +0.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.20"
$ws.Range("E20").Value = "  +2.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.895"
$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("E22").Value = "  +1.24%  "

$ws.Range("E23").Value = "  -0.78%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "100.02"
$ws.Range("E24").Value = "  +1.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.91"
$ws.Range("E25").Value = "  +1.67%  "

$ws.Range("E26").Value = "  +1.29%  "

$ws.Range("E27").Value = "  +2.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "32.99"
$ws.Range("E28").Value = "  -1.83%  "

$ws.Range("E29").Value = "  +1.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.97"
$ws.Range("E30").Value = "  -4.84%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.09"
$ws.Range("E31").Value = "  +1.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "552.35"
$ws.Range("E32").Value = "  -3.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.106"
$ws.Range("E33").Value = "  +0.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "58.17"
$ws.Range("E34").Value = "  +2.69%  "

$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "3.718.84"
$ws.Range("E35").Value = "  +0.68%  "

$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.37"
$ws.Range("E37").Value = "  +0.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").Value = "  +8.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.77"
$ws.Range("E39").Value = "  +1.07%  "

$ws.Range("E40").Value = "  +2.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.61"
$ws.Range("E41").Value = "  -0.26%  "

$ws.Range("E42").Value = "  +1.23%  "

$ws.Range("E43").Value = "  +0.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.27"
$ws.Range("E44").Value = "  -1.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0411"
$ws.Range("E45").Value = "  +1.53%  "

$ws.Range("E46").Value = "  +1.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.128"
$ws.Range("E47").Value = "  +0.36%  "

$ws.Range("E48").Value = "  -0.20%  "

$ws.Range("E49").Value = "  -0.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.64"

$ws.Range("E51").Value = "  -1.15%  "
